$wb = $excel.ActiveWorkbook

$wsDefs = $wb.Worksheets.Item("Defs")
$wsDefs.Range("A2").Value = "lgs5e3pj-elk6"
$wsDefs.Range("B2").Value = "2023-04-22T15:41:12.391Z"
$wsDefs.Range("C2").Value = "lgs5e3pj"

$wsDefs.Range("A3").Value = "lgs5e3pj-6tewf"
$wsDefs.Range("B3").Value = "2023-04-22T15:41:12.391Z"
$wsDefs.Range("C3").Value = "lgs5e3pj"

$wsDefs.Range("A4").Value = "lgs5e3pj-1ve7"
$wsDefs.Range("B4").Value = "2023-04-22T15:41:12.391Z"
$wsDefs.Range("C4").Value = "lgs5e3pj"

$wsPointDefs = $wb.Worksheets.Item("Point Defs")
$wsPointDefs.Range("A2").Value = "lgs5e3pk-0cjl"
$wsPointDefs.Range("B2").Value = "2023-04-22T15:41:12.392Z"
$wsPointDefs.Range("C2").Value = "lgs5e3pk"
$wsPointDefs.Range("F2").Value = "0eze"

$wsPointDefs.Range("A3").Value = "lgs5e3pk-d89q"
$wsPointDefs.Range("B3").Value = "2023-04-22T15:41:12.392Z"
$wsPointDefs.Range("C3").Value = "lgs5e3pk"

$wsPointDefs.Range("A4").Value = "lgs5e3pk-0kt2"
$wsPointDefs.Range("B4").Value = "2023-04-22T15:41:12.392Z"
$wsPointDefs.Range("C4").Value = "lgs5e3pk"

$wsEntryBase = $wb.Worksheets.Item("Entry Base")
$wsEntryBase.Range("A2").Value = "lgs5e3pk-ust9"
$wsEntryBase.Range("B2").Value = "2023-04-22T15:41:12.392Z"
$wsEntryBase.Range("C2").Value = "lgs5e3pk"

$wsEntryBase.Range("A3").Value = "lgs5e3pk-euus"
$wsEntryBase.Range("B3").Value = "2023-04-22T15:41:12.392Z"
$wsEntryBase.Range("C3").Value = "lgs5e3pk"
$wsEntryBase.Range("F3").Value = "lgs5e3pv-5ph5n"
$wsEntryBase.Range("G3").Value = "2023-04-22T10:41:12"

$wsEntryPoints = $wb.Worksheets.Item("Entry Points")
$wsEntryPoints.Range("A2").Value = "lgs5e3pk-00bg"
$wsEntryPoints.Range("B2").Value = "2023-04-22T15:41:12.392Z"
$wsEntryPoints.Range("C2").Value = "lgs5e3pk"

$wsEntryPoints.Range("A3").Value = "lgs5e3pk-5gq2"
$wsEntryPoints.Range("B3").Value = "2023-04-22T15:41:12.392Z"
$wsEntryPoints.Range("C3").Value = "lgs5e3pk"
